$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the D (Price) column as Text first so numeric-looking values
# (e.g. "1.028") are not auto-converted to numbers by Excel, then clear
# the formatting afterwards so the cells keep their original (no explicit
# style) appearance, matching the source workbook.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# Update Price (D) values
$ws.Range("D2").Value = '27.560.17'
$ws.Range("D3").Value = '1.849.21'
$ws.Range("D4").Value = '1.028'
$ws.Range("D5").Value = '320.76'
$ws.Range("D6").Value = '1.027'
$ws.Range("D7").Value = '0.4379'
$ws.Range("D8").Value = '0.3773'
$ws.Range("D9").Value = '0.07386'
$ws.Range("D10").Value = '0.8748'
$ws.Range("D11").Value = '21.49'
$ws.Range("D12").Value = '1.850.58'
$ws.Range("D13").Value = '5.515'
$ws.Range("D14").Value = '6.681'
$ws.Range("D15").Value = '0.07176'
$ws.Range("D16").Value = '82.87'
$ws.Range("D17").Value = '1.033'
$ws.Range("D18").Value = '0.000009014'
$ws.Range("D20").Value = '15.40'
$ws.Range("D21").Value = '27.576.78'
$ws.Range("D22").Value = '5.252'
$ws.Range("D23").Value = '11.35'
$ws.Range("D24").Value = '157.64'
$ws.Range("D25").Value = '1.915'
$ws.Range("D26").Value = '18.71'
$ws.Range("D27").Value = '1.981'
$ws.Range("D28").Value = '5.270'
$ws.Range("D29").Value = '117.05'
$ws.Range("D30").Value = '0.09047'
$ws.Range("D31").Value = '1.195'
$ws.Range("D32").Value = '0.7606'
$ws.Range("D33").Value = '4.521'
$ws.Range("D34").Value = '2.871'
$ws.Range("D37").Value = '0.01976'
$ws.Range("D38").Value = '0.05296'
$ws.Range("D39").Value = '0.5154'
$ws.Range("D40").Value = '2.796'
$ws.Range("D41").Value = '0.1673'
$ws.Range("D42").Value = '6.748'
$ws.Range("D43").Value = '8.482'
$ws.Range("D44").Value = '108.60'
$ws.Range("D45").Value = '10.58'
$ws.Range("D46").Value = '1.707'
$ws.Range("D47").Value = '0.4640'
$ws.Range("D49").Value = '1.855'
$ws.Range("D50").Value = '39.11'
$ws.Range("D51").Value = '63.86'

# Remove the temporary Text formatting so the cells return to the default style
$priceRange.ClearFormats()

# Update Volume(1h) (E) values
$ws.Range("E2").Value = '  +2.45%  '
$ws.Range("E3").Value = '  +2.04%  '
$ws.Range("E4").Value = '  +2.41%  '
$ws.Range("E5").Value = '  +3.11%  '
$ws.Range("E6").Value = '  +2.31%  '
$ws.Range("E7").Value = '  +2.22%  '
$ws.Range("E8").Value = '  +2.47%  '
$ws.Range("E9").Value = '  +2.13%  '
$ws.Range("E10").Value = '  +1.57%  '
$ws.Range("E11").Value = '  +0.82%  '
$ws.Range("E12").Value = '  -8.34%  '
$ws.Range("E13").Value = '  +2.42%  '
$ws.Range("E14").Value = '  +0.86%  '
$ws.Range("E15").Value = '  +4.08%  '
$ws.Range("E16").Value = '  +2.77%  '
$ws.Range("E17").Value = '  +2.46%  '
$ws.Range("E18").Value = '  +0.99%  '
$ws.Range("E19").Value = '  +2.30%  '
$ws.Range("E20").Value = '  +1.00%  '
$ws.Range("E21").Value = '  +2.38%  '
$ws.Range("E22").Value = '  +1.41%  '
$ws.Range("E23").Value = '  +3.13%  '
$ws.Range("E24").Value = '  +2.57%  '
$ws.Range("E25").Value = '  +1.69%  '
$ws.Range("E26").Value = '  +2.42%  '
$ws.Range("E27").Value = '  +5.32%  '
$ws.Range("E28").Value = '  +0.98%  '
$ws.Range("E29").Value = '  +1.86%  '
$ws.Range("E30").Value = '  +1.13%  '
$ws.Range("E31").Value = '  +2.70%  '
$ws.Range("E32").Value = '  +2.42%  '
$ws.Range("E33").Value = '  +2.09%  '
$ws.Range("E34").Value = '  +2.62%  '
$ws.Range("E35").Value = '  +1.93%  '
$ws.Range("E36").Value = '  +2.92%  '
$ws.Range("E37").Value = '  +2.85%  '
$ws.Range("E38").Value = '  +1.52%  '
$ws.Range("E39").Value = '  +1.38%  '
$ws.Range("E40").Value = '  +2.50%  '
$ws.Range("E41").Value = '  +1.79%  '
$ws.Range("E42").Value = '  +4.77%  '
$ws.Range("E43").Value = '  +2.84%  '
$ws.Range("E44").Value = '  +1.67%  '
$ws.Range("E45").Value = '  +2.12%  '
$ws.Range("E46").Value = '  +3.27%  '
$ws.Range("E47").Value = '  +1.32%  '
$ws.Range("E48").Value = '  +1.75%  '
$ws.Range("E49").Value = '  +2.64%  '
$ws.Range("E50").Value = '  +3.94%  '
$ws.Range("E51").Value = '  +0.35%  '
